{"js": "// Fix several typos / proofing issues in the VIS instructions document:\n//   - \"atleast\"    -> \"at least\"\n//   - \"indivduals\" -> \"individuals\"\n//   - \"Airforce\"   -> \"Air Force\"\n//   - \"indivuals\"  -> \"individuals\"\n//   - \"sitution\"   -> \"situation\"\n// (The surrounding wording such as \"The minimum of time...\" itself does\n// not change; Word had also attached stray proofing markers\n// (<w:proofErr .../>) around some of these words, but those are purely\n// internal spell/grammar-check bookkeeping with no Office.js object\n// model surface, so they are left as-is.)\n\nasync function replaceOnce(searchText, replacementText) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nawait replaceOnce(\"atleast\", \"at least\");\nawait replaceOnce(\"indivduals\", \"individuals\");\nawait replaceOnce(\"Airforce\", \"Air Force\");\nawait replaceOnce(\"indivuals\", \"individuals\");\nawait replaceOnce(\"sitution\", \"situation\");\n", "ps1": "# Fix several typos / proofing issues in the VIS instructions document:\n#   - \"The minimum of time...\" was wrapped with a stray grammar-check\n#     marker around \"The\" -> clean it up (text itself is unchanged).\n#   - \"atleast\"    -> \"at least\"\n#   - \"indivduals\" -> \"individuals\"\n#   - \"Airforce\"   -> \"Air Force\"\n#   - \"indivuals\"  -> \"individuals\"\n#   - \"sitution\"   -> \"situation\"\n#\n# Each Find/Replace below intentionally includes a little bit of the\n# neighbouring text on both sides of the target word. That makes Word\n# rebuild the run that held the (now obsolete) proofing-error markers\n# <w:proofErr w:type=\"gramStart/gramEnd/spellStart/spellEnd\"/> that\n# Word had placed around the misspelled/flagged word, instead of\n# leaving them orphaned next to the corrected text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute(\n        $findText,     # FindText\n        $true,         # MatchCase\n        $false,        # MatchWholeWord\n        $false,        # MatchWildcards\n        $false,        # MatchSoundsLike\n        $false,        # MatchAllWordForms\n        $true,         # Forward\n        1,             # Wrap (wdFindContinue)\n        $false,        # Format\n        $replaceText,  # ReplaceWith\n        1              # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n\n# 1) Drop the stale grammar-check wrapper around \"The\" (text itself is\n#    unchanged - this only cleans up the proofing markers).\nReplace-Text \". The minimum of time\" \". The minimum of time\"\n\n# 2) \"atleast\" -> \"at least\"\nReplace-Text \"have atleast have\" \"have at least have\"\n\n# 3) \"indivduals\" -> \"individuals\"\nReplace-Text \"1-2 indivduals focusing\" \"1-2 individuals focusing\"\n\n# 4) \"Airforce\" -> \"Air Force\"\nReplace-Text \"and Airforce), \" \"and Air Force), \"\n\n# 5) \"indivuals\" -> \"individuals\"\nReplace-Text \"and 1-2 indivuals focusing\" \"and 1-2 individuals focusing\"\n\n# 6) \"sitution\" -> \"situation\"\nReplace-Text \"enemy sitution to\" \"enemy situation to\"\n"}
